$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 33 entirely; this shifts rows 34:56 up to 33:55
$ws.Rows("33:33").Delete()

# Update a handful of C/D values that changed as part of this commit
# (post-deletion row numbers)
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(8, 4).Value = 0.0008324633664117887

$ws.Cells.Item(24, 3).Value = 1
$ws.Cells.Item(24, 4).Value = 0.005786472390968115

$ws.Cells.Item(25, 3).Value = 0
$ws.Cells.Item(25, 4).Value = 0.003175346061097805

$ws.Cells.Item(30, 3).Value = 4
$ws.Cells.Item(30, 4).Value = 0.0507632814884952

$ws.Cells.Item(32, 3).Value = 2
$ws.Cells.Item(32, 4).Value = 0.04066896040262172

# Column A holds a literal running index (0,1,2,...), not a formula, so after
# deleting row 33 the index values below it must be renumbered down by one to
# stay sequential (31, 32, 33, ... instead of 32, 33, 34, ...).
for ($r = 33; $r -le 55; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
